$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) values from 45175 to 45183 for rows 2-27
$ws.Range("C2:C27").Value = 45183
